# The two ~TFM_INS blocks (rows 3-12 and rows 16-25) had their data swapped:
# rows 3-12 (previously a run of ACT_COST / TU_* fuel rows) now carry the
# VAROM/LIFE/CAP2ACT/... TB_ELC* attribute rows that used to live in rows
# 16-25, and vice versa. Row 10 / row 23 additionally carry a Year value
# ("2010") for the NCAP_PASTI attribute, which moves along with the rest of
# that row's data.
#
# Because the source cells store every value (including numeric-looking
# ones) as literal text (inline strings), we force each written cell to
# Text format before assigning it, then clear that cell's formatting again
# so we don't leave a stray NumberFormat behind (the source file uses no
# per-cell styles at all).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.ClearFormats()
}

# New content for rows 3-12 (column C / D / F / H). $null means "no D value".
$group1 = @{
    3  = @("VAROM",      $null,   "0",       "TB_ELC*")
    4  = @("LIFE",       $null,   "100",     "TB_ELC*")
    5  = @("CAP2ACT",    $null,   "31.536",  "TB_ELC*")
    6  = @("PEAK(CON)",  $null,   "1",       "TB_ELC*")
    7  = @("AFA",        $null,   "1",       "TB_ELC*")
    8  = @("EFF",        $null,   "0.99",    "TB_ELC*")
    9  = @("NCAP_ILED",  $null,   "3",       "TB_ELC*")
    10 = @("NCAP_PASTI", "2010",  "1.6",     "TB_ELC_*")
    11 = @("INVCOST",    $null,   "480.0",   "TB_ELC*")
    12 = @("FIXOM",      $null,   "0.96",    "TB_ELC*")
}

foreach ($r in $group1.Keys) {
    $vals = $group1[$r]
    Set-TextValue $r 3 $vals[0]
    if ($vals[1] -ne $null) {
        Set-TextValue $r 4 $vals[1]
    }
    Set-TextValue $r 6 $vals[2]
    Set-TextValue $r 8 $vals[3]
}

# New content for rows 16-25 (column C / F / H) - all become ACT_COST rows.
$group2 = @{
    16 = @("ACT_COST", "0.00248362891109477", "TU_PET*")
    17 = @("ACT_COST", "0.00538337115666179", "TU_LPG*")
    18 = @("ACT_COST", "0.00252841765861999", "TU_DSL*")
    19 = @("ACT_COST", "0.00261893423255859", "TU_FOL*")
    20 = @("ACT_COST", "0.00252841765861999", "TU_DID*")
    21 = @("ACT_COST", "0.00252841765861999", "TU_DIJ*")
    22 = @("ACT_COST", "0.00251030880245087", "TU_JET*")
    23 = @("ACT_COST", "0.0025",              "TU_OTH*")
    24 = @("ACT_COST", "0.000188040616773223","TU_COA*")
    25 = @("ACT_COST", "0.000327439423706614","TU_COL*")
}

foreach ($r in $group2.Keys) {
    $vals = $group2[$r]
    Set-TextValue $r 3 $vals[0]
    Set-TextValue $r 6 $vals[1]
    Set-TextValue $r 8 $vals[2]
}

# Row 23's Year cell (D23) no longer holds "2010" - it goes back to being an
# empty numeric placeholder cell, matching the other rows in that block.
$ws.Cells.Item(23, 4).Value2 = 0
